# Weekly price update: insert two new rows (a new reporting week) right
# after the current most-recent row (row 349), pushing all the
# historical rows down by two. The new rows carry the same constant
# columns (Mercado / Region / Categoria / Variedad / Unidad / Origen /
# Kg-o-Unidades / Clasificacion) as the rest of the sheet, plus the new
# date and price figures for quality grades "Segunda" and "Tercera".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 350 - everything that was on row 350
# onward (through 382) shifts down to 352..384.
$ws.Range("A350:A351").EntireRow.Insert()

# --- Row 350 : Segunda ---
$ws.Range("A350").Value = 1
$ws.Range("B350").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C350").Value = "Arica y Parinacota"
$ws.Range("D350").Value = 44714
$ws.Range("D350").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E350").Value = 15
$ws.Range("F350").Value = 100112023
$ws.Range("G350").Value = "Brócoli"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Segunda"
$ws.Range("J350").Value = 700
$ws.Range("K350").Value = 600
$ws.Range("L350").Value = 700
$ws.Range("M350").Value = 650
$ws.Range("N350").Value = "$/unidad"
$ws.Range("O350").Value = "Región de Arica y Parinacota"
$ws.Range("P350").Value = 650
$ws.Range("Q350").Value = 1
$ws.Range("R350").Value = "Hortaliza"

# --- Row 351 : Tercera ---
$ws.Range("A351").Value = 1
$ws.Range("B351").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C351").Value = "Arica y Parinacota"
$ws.Range("D351").Value = 44714
$ws.Range("D351").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E351").Value = 15
$ws.Range("F351").Value = 100112023
$ws.Range("G351").Value = "Brócoli"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Tercera"
$ws.Range("J351").Value = 1200
$ws.Range("K351").Value = 400
$ws.Range("L351").Value = 500
$ws.Range("M351").Value = 450
$ws.Range("N351").Value = "$/unidad"
$ws.Range("O351").Value = "Región de Arica y Parinacota"
$ws.Range("P351").Value = 450
$ws.Range("Q351").Value = 1
$ws.Range("R351").Value = "Hortaliza"
